# Change the "动态寻址（可选）" header/switcher into a "寻址方式（可选）" select
# with options "静态" / "动态", per commit: 寻址方式 switcher 改为 select

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in O1 (was "动态寻址（可选）") to new label "寻址方式（可选）"
$ws.Range("O1").Value = "寻址方式（可选）"

# Row 2 used to hold a boolean TRUE switch value; now holds the selected option text "静态"
$ws.Range("O2").Value = "静态"

# Row 3 gains a new value for the alternate option "动态"
$ws.Range("O3").Value = "动态"

# Update the active selection to O3 to match the saved view state
$ws.Range("O3").Select()
